# Adds a new weekly price record for "Feria Lagunitas de Puerto Montt - Espinaca".
# The new record is inserted as row 17 (date 2022-08-26 / serial 44799), pushing
# the existing rows 17-38 down to rows 18-39 and extending the used range to
# A1:R39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 17, shifting rows 17:38 down to 18:39.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row with the new market observation.
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44799
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112012
$ws.Range("G17").Value = "Espinaca"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("N17").Value = "$/cuna 10 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 1200
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = "Hortaliza"
